$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 data first (so new shared strings are interned in this order)
$ws.Range("A6").Value = "localhost"
$ws.Range("B6").Value = "UMA_DWH"
$ws.Range("C6").Value = "MWH"
$ws.Range("D6").Value = "AgentActivityLog"
$ws.Range("E6").Value = "MARKETING"
$ws.Range("F6").Value = "I3TimeStampGMT"
$ws.Range("G6").Value = $true
$ws.Range("K6").Value = "localhost"
$ws.Range("L6").Value = "UMA_DWH"
$ws.Range("N6").Value = "AgentActivityLog"
$ws.Range("Q6").Value = "PRIMARY"
$ws.Range("S6").Value = "PRIMARY"
$ws.Range("U6").Value = "MERGE_<TARGET_TABLE>"

# Update column M (TARGET_SCHEMA) values on rows 2-5 from "MWH" to "TEST"
$ws.Range("M2:M5").Value = "TEST"
$ws.Range("M6").Value = "TEST"

# Update the view/selection
$ws.Range("V7").Select()
